# Refresh the "Neg_Change" and "Pos_Change" market-data sheets with the
# latest scraped rows (GitHub Actions data refresh).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Neg_Change")
$ws.Cells.Item(2,1).Value = "HINDUNILVR"
$ws.Cells.Item(2,2).Value = 2302.6
$ws.Cells.Item(2,3).Value = 2305
$ws.Cells.Item(2,4).Value = 2278.1
$ws.Cells.Item(2,5).Value = 2283
$ws.Cells.Item(2,6).Value = 519380
$ws.Cells.Item(2,7).Value = 1128801
$ws.Cells.Item(2,8).Value = -0.539883469274035
$ws.Cells.Item(2,9).Value = "HINDUNILVR"
$ws.Cells.Item(3,1).Value = "RECLTD"
$ws.Cells.Item(3,2).Value = 355.1
$ws.Cells.Item(3,3).Value = 358.2
$ws.Cells.Item(3,4).Value = 353.8
$ws.Cells.Item(3,5).Value = 354.9
$ws.Cells.Item(3,6).Value = 3790465
$ws.Cells.Item(3,7).Value = 8591989
$ws.Cells.Item(3,8).Value = -0.5588373076362179
$ws.Cells.Item(3,9).Value = "RECLTD"
$ws.Cells.Item(4,1).Value = "IRFC"
$ws.Cells.Item(4,2).Value = 122.5
$ws.Cells.Item(4,3).Value = 123.84
$ws.Cells.Item(4,4).Value = 120.54
$ws.Cells.Item(4,5).Value = 121.54
$ws.Cells.Item(4,6).Value = 27580407
$ws.Cells.Item(4,7).Value = 55006061
$ws.Cells.Item(4,8).Value = -0.4985933095627407
$ws.Cells.Item(4,9).Value = "IRFC"
$ws.Cells.Item(5,1).Value = "LICI"
$ws.Cells.Item(5,2).Value = 855.4
$ws.Cells.Item(5,3).Value = 858
$ws.Cells.Item(5,4).Value = 851.5
$ws.Cells.Item(5,5).Value = 851.7
$ws.Cells.Item(5,6).Value = 494054
$ws.Cells.Item(5,7).Value = 1143547
$ws.Cells.Item(5,8).Value = -0.5679635380093691
$ws.Cells.Item(5,9).Value = "LICI"
$ws.Cells.Item(6,1).Value = "BRITANNIA"
$ws.Cells.Item(6,2).Value = 6065
$ws.Cells.Item(6,3).Value = 6078
$ws.Cells.Item(6,4).Value = 6004
$ws.Cells.Item(6,5).Value = 6016
$ws.Cells.Item(6,6).Value = 97580
$ws.Cells.Item(6,7).Value = 216840
$ws.Cells.Item(6,8).Value = -0.5499907766094816
$ws.Cells.Item(6,9).Value = "BRITANNIA"
$ws.Cells.Item(7,1).Value = "PFC"
$ws.Cells.Item(7,2).Value = 356
$ws.Cells.Item(7,3).Value = 357
$ws.Cells.Item(7,4).Value = 351.05
$ws.Cells.Item(7,5).Value = 352
$ws.Cells.Item(7,6).Value = 4097578
$ws.Cells.Item(7,7).Value = 9432271
$ws.Cells.Item(7,8).Value = -0.5655788515830387
$ws.Cells.Item(7,9).Value = "PFC"
$ws.Cells.Item(8,1).Value = "OIL"
$ws.Cells.Item(8,2).Value = 410.4
$ws.Cells.Item(8,3).Value = 413.5
$ws.Cells.Item(8,4).Value = 408
$ws.Cells.Item(8,5).Value = 408.9
$ws.Cells.Item(8,6).Value = 706490
$ws.Cells.Item(8,7).Value = 1528489
$ws.Cells.Item(8,8).Value = -0.5377853553411245
$ws.Cells.Item(8,9).Value = "OIL"
$ws.Cells.Item(9,1).Value = "SUPREMEIND"
$ws.Cells.Item(9,2).Value = 3360
$ws.Cells.Item(9,3).Value = 3363
$ws.Cells.Item(9,4).Value = 3314.3
$ws.Cells.Item(9,5).Value = 3327.9
$ws.Cells.Item(9,6).Value = 77105
$ws.Cells.Item(9,7).Value = 152437
$ws.Cells.Item(9,8).Value = -0.4941844827699312
$ws.Cells.Item(9,9).Value = "SUPREMEIND"
$ws.Cells.Item(10,1).Value = "CONCOR"
$ws.Cells.Item(10,2).Value = 515.05
$ws.Cells.Item(10,3).Value = 515.3
$ws.Cells.Item(10,4).Value = 506.75
$ws.Cells.Item(10,5).Value = 509.05
$ws.Cells.Item(10,6).Value = 476260
$ws.Cells.Item(10,7).Value = 943869
$ws.Cells.Item(10,8).Value = -0.4954172665910206
$ws.Cells.Item(10,9).Value = "CONCOR"
$ws.Cells.Item(11,1).Value = "PIIND"
$ws.Cells.Item(11,2).Value = 3239.7
$ws.Cells.Item(11,3).Value = 3257.2
$ws.Cells.Item(11,4).Value = 3210.1
$ws.Cells.Item(11,5).Value = 3215
$ws.Cells.Item(11,6).Value = 74555
$ws.Cells.Item(11,7).Value = 148650
$ws.Cells.Item(11,8).Value = -0.4984527413387151
$ws.Cells.Item(11,9).Value = "PIIND"
$ws.Cells.Item(12,1).Value = "PATANJALI"
$ws.Cells.Item(12,2).Value = 552.05
$ws.Cells.Item(12,3).Value = 554.7
$ws.Cells.Item(12,4).Value = 540.65
$ws.Cells.Item(12,5).Value = 545.9
$ws.Cells.Item(12,6).Value = 1918838
$ws.Cells.Item(12,7).Value = 3784328
$ws.Cells.Item(12,8).Value = -0.4929514566390651
$ws.Cells.Item(12,9).Value = "PATANJALI"
$ws.Cells.Item(13,1).Value = "GLENMARK"
$ws.Cells.Item(13,2).Value = 2050
$ws.Cells.Item(13,3).Value = 2054.9
$ws.Cells.Item(13,4).Value = 2013.8
$ws.Cells.Item(13,5).Value = 2015.4
$ws.Cells.Item(13,6).Value = 579266
$ws.Cells.Item(13,7).Value = 1403145
$ws.Cells.Item(13,8).Value = -0.5871659735807775
$ws.Cells.Item(13,9).Value = "GLENMARK"

$ws = $wb.Worksheets.Item("Pos_Change")
$ws.Cells.Item(2,1).Value = "TRENT"
$ws.Cells.Item(2,2).Value = 4209
$ws.Cells.Item(2,3).Value = 4313.3
$ws.Cells.Item(2,4).Value = 4183.1
$ws.Cells.Item(2,5).Value = 4284
$ws.Cells.Item(2,6).Value = 960632
$ws.Cells.Item(2,7).Value = 662345
$ws.Cells.Item(2,8).Value = 0.4503498931825559
$ws.Cells.Item(2,9).Value = "TRENT"
$ws.Cells.Item(3,1).Value = "ADANIPORTS"
$ws.Cells.Item(3,2).Value = 1503
$ws.Cells.Item(3,3).Value = 1507.2
$ws.Cells.Item(3,4).Value = 1492
$ws.Cells.Item(3,5).Value = 1496.5
$ws.Cells.Item(3,6).Value = 1318023
$ws.Cells.Item(3,7).Value = 928446
$ws.Cells.Item(3,8).Value = 0.4196011399693682
$ws.Cells.Item(3,9).Value = "ADANIPORTS"
$ws.Cells.Item(4,1).Value = "HINDZINC"
$ws.Cells.Item(4,2).Value = 619
$ws.Cells.Item(4,3).Value = 631.9
$ws.Cells.Item(4,4).Value = 617.95
$ws.Cells.Item(4,5).Value = 625.45
$ws.Cells.Item(4,6).Value = 13964315
$ws.Cells.Item(4,7).Value = 9277131
$ws.Cells.Item(4,8).Value = 0.5052406827067549
$ws.Cells.Item(4,9).Value = "HINDZINC"
$ws.Cells.Item(5,1).Value = "ABB"
$ws.Cells.Item(5,2).Value = 5212.5
$ws.Cells.Item(5,3).Value = 5266
$ws.Cells.Item(5,4).Value = 5202
$ws.Cells.Item(5,5).Value = 5202.5
$ws.Cells.Item(5,6).Value = 103441
$ws.Cells.Item(5,7).Value = 64953
$ws.Cells.Item(5,8).Value = 0.5925515372654073
$ws.Cells.Item(5,9).Value = "ABB"
$ws.Cells.Item(6,1).Value = "JSWENERGY"
$ws.Cells.Item(6,2).Value = 486.9
$ws.Cells.Item(6,3).Value = 488.1
$ws.Cells.Item(6,4).Value = 481.55
$ws.Cells.Item(6,5).Value = 483.55
$ws.Cells.Item(6,6).Value = 1822866
$ws.Cells.Item(6,7).Value = 1158705
$ws.Cells.Item(6,8).Value = 0.5731924864396029
$ws.Cells.Item(6,9).Value = "JSWENERGY"
$ws.Cells.Item(7,1).Value = "BOSCHLTD"
$ws.Cells.Item(7,2).Value = 36210
$ws.Cells.Item(7,3).Value = 36470
$ws.Cells.Item(7,4).Value = 36010
$ws.Cells.Item(7,5).Value = 36150
$ws.Cells.Item(7,6).Value = 16626
$ws.Cells.Item(7,7).Value = 11103
$ws.Cells.Item(7,8).Value = 0.4974331261821129
$ws.Cells.Item(7,9).Value = "BOSCHLTD"
$ws.Cells.Item(8,1).Value = "DMART"
$ws.Cells.Item(8,2).Value = 3843
$ws.Cells.Item(8,3).Value = 3843
$ws.Cells.Item(8,4).Value = 3771.3
$ws.Cells.Item(8,5).Value = 3805
$ws.Cells.Item(8,6).Value = 363697
$ws.Cells.Item(8,7).Value = 240671
$ws.Cells.Item(8,8).Value = 0.5111791615940433
$ws.Cells.Item(8,9).Value = "DMART"
$ws.Cells.Item(9,1).Value = "CGPOWER"
$ws.Cells.Item(9,2).Value = 666.3
$ws.Cells.Item(9,3).Value = 669.25
$ws.Cells.Item(9,4).Value = 657.15
$ws.Cells.Item(9,5).Value = 660.3
$ws.Cells.Item(9,6).Value = 1910047
$ws.Cells.Item(9,7).Value = 1294033
$ws.Cells.Item(9,8).Value = 0.476041955653372
$ws.Cells.Item(9,9).Value = "CGPOWER"
$ws.Cells.Item(10,1).Value = "SUZLON"
$ws.Cells.Item(10,2).Value = 53.59
$ws.Cells.Item(10,3).Value = 55.05
$ws.Cells.Item(10,4).Value = 53.16
$ws.Cells.Item(10,5).Value = 53.3
$ws.Cells.Item(10,6).Value = 44076962
$ws.Cells.Item(10,7).Value = 30259539
$ws.Cells.Item(10,8).Value = 0.4566303207725669
$ws.Cells.Item(10,9).Value = "SUZLON"
$ws.Cells.Item(11,1).Value = "ALKEM"
$ws.Cells.Item(11,2).Value = 5570
$ws.Cells.Item(11,3).Value = 5593
$ws.Cells.Item(11,4).Value = 5534
$ws.Cells.Item(11,5).Value = 5552
$ws.Cells.Item(11,6).Value = 41257
$ws.Cells.Item(11,7).Value = 28301
$ws.Cells.Item(11,8).Value = 0.4577930108476732
$ws.Cells.Item(11,9).Value = "ALKEM"
$ws.Cells.Item(12,1).Value = "BIOCON"
$ws.Cells.Item(12,2).Value = 404
$ws.Cells.Item(12,3).Value = 406.35
$ws.Cells.Item(12,4).Value = 397.35
$ws.Cells.Item(12,5).Value = 398.05
$ws.Cells.Item(12,6).Value = 2707789
$ws.Cells.Item(12,7).Value = 1907267
$ws.Cells.Item(12,8).Value = 0.4197220420633294
$ws.Cells.Item(12,9).Value = "BIOCON"
$ws.Cells.Item(13,1).Value = "INDIANB"
$ws.Cells.Item(13,2).Value = 783.25
$ws.Cells.Item(13,3).Value = 783.25
$ws.Cells.Item(13,4).Value = 773.25
$ws.Cells.Item(13,5).Value = 774.65
$ws.Cells.Item(13,6).Value = 902942
$ws.Cells.Item(13,7).Value = 581636
$ws.Cells.Item(13,8).Value = 0.5524176632808148
$ws.Cells.Item(13,9).Value = "INDIANB"
$ws.Cells.Item(14,1).Value = "INOXWIND"
$ws.Cells.Item(14,2).Value = 127.1
$ws.Cells.Item(14,3).Value = 130.9
$ws.Cells.Item(14,4).Value = 126.4
$ws.Cells.Item(14,5).Value = 126.7
$ws.Cells.Item(14,6).Value = 8410588
$ws.Cells.Item(14,7).Value = 5458544
$ws.Cells.Item(14,8).Value = 0.5408116156982521
$ws.Cells.Item(14,9).Value = "INOXWIND"
$ws.Cells.Item(15,1).Value = "BANDHANBNK"
$ws.Cells.Item(15,2).Value = 149.9
$ws.Cells.Item(15,3).Value = 149.9
$ws.Cells.Item(15,4).Value = 146.13
$ws.Cells.Item(15,5).Value = 146.35
$ws.Cells.Item(15,6).Value = 3955280
$ws.Cells.Item(15,7).Value = 2568295
$ws.Cells.Item(15,8).Value = 0.5400411557083591
$ws.Cells.Item(15,9).Value = "BANDHANBNK"
$ws.Cells.Item(16,1).Value = "IEX"
$ws.Cells.Item(16,2).Value = 142.94
$ws.Cells.Item(16,3).Value = 143.4
$ws.Cells.Item(16,4).Value = 138.8
$ws.Cells.Item(16,5).Value = 139.13
$ws.Cells.Item(16,6).Value = 3752961
$ws.Cells.Item(16,7).Value = 2545162
$ws.Cells.Item(16,8).Value = 0.4745470032948787
$ws.Cells.Item(16,9).Value = "IEX"
